# Fix batch ingest of single-value fields
#
# Adds two new "single value" metadata columns — Abstract and Statement Of
# Responsibility — to the example batch-ingest manifest, along with sample
# values for the two example rows used by the fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the added columns (row 2 holds the field-name headers).
$ws.Range("Y2").Value = "Abstract"
$ws.Range("Z2").Value = "Statement Of Responsibility"

# First example item (row 3) gets both single-value fields populated.
$ws.Range("Y3").Value = "Test abstract"
$ws.Range("Z3").Value = "Test Statement of Responsibility"

# Second example item (row 4) only has the Abstract populated.
$ws.Range("Y4").Value = "Test abstract"

# Match the workbook's saved selection/view state.
$ws.Range("Y4").Select()
